$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Disgust" row (row 6) entirely, shifting "Surprise" (row 7) up to row 6.
$ws.Rows.Item(6).Delete()

# Update header row text
$ws.Range("B1").Value = "Hume (speech)"
$ws.Range("C1").Value = "NLP (text)"
$ws.Range("D1").Value = "Self$([char]0x2011)label"

# Copy the header style (bold, centered, thin border) from an existing header cell onto D1
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Updated data values for columns B and C, plus new column D
$ws.Range("B2").Value = 0.26
$ws.Range("C2").Value = 0.22
$ws.Range("D2").Value = 0.08

$ws.Range("B3").Value = 0.15
$ws.Range("C3").Value = 0.35
$ws.Range("D3").Value = 0.38

$ws.Range("B4").Value = 0.17
$ws.Range("C4").Value = 0.09
$ws.Range("D4").Value = 0.08

$ws.Range("B5").Value = 0.31
$ws.Range("C5").Value = 0.04
$ws.Range("D5").Value = 0.08

$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 0.3
$ws.Range("D6").Value = 0.38
